# Updates odds values on rows 20-22 of Sheet1 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20
$ws.Range("G20").Value = 1.73
$ws.Range("I20").Value = 4.2
$ws.Range("L20").Value = 4.5
$ws.Range("U20").Value = 1.67
$ws.Range("V20").Value = 2.1
$ws.Range("W20").Value = 8.5
$ws.Range("Z20").Value = 15
$ws.Range("AG20").Value = 151
$ws.Range("AH20").Value = 13
$ws.Range("AK20").Value = 41
$ws.Range("AN20").Value = 4
$ws.Range("AO20").Value = 9
$ws.Range("AQ20").Value = 29
$ws.Range("AU20").Value = 7.5
$ws.Range("AY20").Value = 26
$ws.Range("AZ20").Value = 67

# Row 21
$ws.Range("O21").Value = 1.17
$ws.Range("P21").Value = 5

# Row 22
$ws.Range("G22").Value = 1.36
$ws.Range("J22").Value = 1.8
$ws.Range("L22").Value = 6
$ws.Range("X22").Value = 9.5
$ws.Range("Z22").Value = 11
$ws.Range("AD22").Value = 11
$ws.Range("AX22").Value = 29
$ws.Range("BA22").Value = 81
$ws.Range("BB22").Value = 126
